$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '$2b$10$sZ9I.kZKPGy.I/KTpI5T8.AFs/is1186iazck4NjzqaAaUGdCXUeS'
$ws.Range("C3").Value = '$2b$10$0SehfwimNn99Ry4wx67EK.BlEEUTvXK0towkfF7pBlJj/rliymALy'
$ws.Range("C4").Value = '$2b$10$oPS4zDbAsMk34rnjfQwZleN/isd91acC58FW3vYrHidfQbTukNmlG'
$ws.Range("C5").Value = '$2b$10$/cFPBjkc3UetzO3KUP03S..MMVjjEz2Rm6Yqre/l/PK4El.8J6JHC'
$ws.Range("C6").Value = '$2b$10$55ByGzitYQHrru.8AnjSxuw1zPS.C27RtM5hbFJNLsT09VuuofIke'
$ws.Range("C7").Value = '$2b$10$YPXcA20ejcwBpKg1gfGEeuXwrRIyKF9X4eyydkI6ZU57/Y2mZALOO'
$ws.Range("C8").Value = '$2b$10$irTvl2EQ6wlG2Wv9AYKb5uuOaQAs5W8rwcQFXUZnRfUWZ9xWCufLW'
$ws.Range("C9").Value = '$2b$10$5RBb8wiTN.smh8Ku3QzVqea.7drhdd080CNVY.NOMhQYoJVb90Epm'
$ws.Range("C10").Value = '$2b$10$VRy8CdwHSIsV6kfG969ujOG9fk/Se2kYp1FhpKp1z1va3s9GaJG7S'
$ws.Range("C11").Value = '$2b$10$NX6o9XKLPoYAvXlmDLeq2eZnSSZkn0yuA2ULjldKqExNQByELkBwq'
$ws.Range("C12").Value = '$2b$10$ZboSh6qZG1IyjobcD/3J3u6ummjcLGudjSkagbTyDVqsLpzb0MP7K'
$ws.Range("C13").Value = '$2b$10$evUKf4g0HIZuN.SqFf7UDOscnmvJn7OXZvaYqn3E1xp5fuI6RcL4G'
$ws.Range("C14").Value = '$2b$10$s6m3iyrtGhn4y/pTinLNjutieRHaMtrdDpuHeb2W0jE4cXo0oj6dG'
$ws.Range("C15").Value = '$2b$10$hxRCnhKBOwFlO/FWrzzYrOQPxzVuNwSe84BZnqataGbt0InoaRht2'
$ws.Range("C16").Value = '$2b$10$s8/pAqjiL7kd2I8rJxaagu49gW8aFcCvoD99ICZpVMal/mEe8v/B2'
$ws.Range("C17").Value = '$2b$10$TaeCGp5/4Jr0SncOiyRCQuoAKIxZ/Xo2cIEE33vPCYjniDF.uZTB.'
$ws.Range("C18").Value = '$2b$10$qVeeMSW3qJKFvC5ziAON0OpkaiadsvBpounkzyRxluJFHxFB83dma'
$ws.Range("C19").Value = '$2b$10$6tQmAirHV92YMV5Jm17reuxAOr05zWh5ZgWDRn1K2B4X24J9ZmRM2'
$ws.Range("C20").Value = '$2b$10$yKx2pp5hQRlyBObOjXcnJejHKH97Mt8p6rhffc9H7YkIL72X9dZCG'
$ws.Range("C21").Value = '$2b$10$8ZCe3Ej86mFfvnvQxrWqO.U2O7eV5lrsBn4VLiG6nzkhQt16b4FMy'
$ws.Range("C22").Value = '$2b$10$F4MBFxVvQDHfcp6ZAE86..3URH5ndxU/VVTbU3QsKLk6NbF.uvzg2'
$ws.Range("C23").Value = '$2b$10$rM/s0OpdpcidLJGhyYRnMermYbrXaN.wGcZ3IsLugD2UMoDPS/Wcu'
$ws.Range("C24").Value = '$2b$10$.2QCHEDqcPFJ/CgT5t7LcO9gge657lH.Vk4BEbxbI1.BtE3EIimrm'
$ws.Range("C25").Value = '$2b$10$p50hTk3v2cPELVAW3kDw2OgyEol6y6EPYUlHmJ6qEScwXIHkN35/u'
$ws.Range("C26").Value = '$2b$10$VdDRL..NFjg1jRvxVcbokOBSPpn2lyRU7KiJ0lNJyGAXfOeOYvdbK'